$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four species rows that were dropped from the dataset
# (GRPC, HOLA, TUVU, AMKE). Deleting shifts everything below up,
# which is exactly the row-shift pattern seen across rows 8-28.
$ws.Rows(8).Delete()   # GRPC (was row 8)
$ws.Rows(22).Delete()  # HOLA (was row 23, now row 22 after the above delete)
$ws.Rows(22).Delete()  # TUVU (was row 24, now row 22)
$ws.Rows(22).Delete()  # AMKE (was row 25, now row 22)

# Reselect the block that was being reviewed/filled in
$ws.Range("D11:P21").Select()

# Add 3-colour-scale conditional formatting (red/yellow/green) to the
# four ranges in the sheet, matching priorities 3, 6, 8, 9
$cf1 = $ws.Range("D22:P28").FormatConditions.AddColorScale(3)
$cf1.Priority = 3

$cf2 = $ws.Range("D2:O10").FormatConditions.AddColorScale(3)
$cf2.Priority = 6

$cf3 = $ws.Range("P2:P10").FormatConditions.AddColorScale(3)
$cf3.Priority = 8

$cf4 = $ws.Range("D11:P21").FormatConditions.AddColorScale(3)
$cf4.Priority = 9

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
